# "Import some images into resource library" — adds a new column E to
# Sheet1 ("图像是否已导入资源库" / "has this image been imported into the
# resource library") with "Y" marked for every row that already has an
# image reference (A/B/C populated in the first, contiguous image block),
# widens column A slightly, and leaves the view scrolled/selected on the
# new column.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rows that get a "Y" flag in the new column E (everything in the first
# "elements" block of the sheet, rows 2-18 plus 35-38).
$yRows = @(2,3,4,5,6,7,8,9,10,11,12,13,14,15,16,17,18,35,36,37,38)

# Write the repeated "Y" values first so that shared string gets created
# (and reused) before the header string, matching the order new strings
# were appended to the workbook.
foreach ($r in $yRows) {
    $ws.Cells.Item($r, 5).Value = "Y"
}

# Header for the new column.
$ws.Range("E1").Value = "图像是否已导入资源库"

# Column A is a little wider now that there's more going on in the sheet.
$ws.Columns.Item(1).ColumnWidth = 18.571428571428571

# Leave the sheet scrolled back to the top with the new column selected
# (this also clears the old topLeftCell/selection that pointed at B80).
$ws.Range("E1:E40").Select()
